# Workers-rank-matrices sheet was regenerated upstream: several ties in the
# "matrices" (C) ranking were broken differently, shifting the prolificid/
# name (and, where applicable, race) pairing between a few adjacent rows,
# and every mat_rank (G) score was recomputed with slightly different
# precision. level_0 / index / gender / rank-position (A, B, F, I) are
# unaffected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Female group (rows 2-13) ---
$ws.Range("G2").Value  = 13.37383182294894
$ws.Range("G3").Value  = 13.16250246588412
$ws.Range("G4").Value  = 8.396910801783761

$ws.Range("C5").Value  = 19
$ws.Range("D5").Value  = "60b45e9961dd412bfb6780f8"
$ws.Range("E5").Value  = "Jewel"
$ws.Range("G5").Value  = 8.390562821666926

$ws.Range("C6").Value  = 21
$ws.Range("D6").Value  = "5c0e89c6c323400001e6c4a5"
$ws.Range("E6").Value  = "Bri"
$ws.Range("G6").Value  = 8.284137808845447

$ws.Range("C7").Value  = 32
$ws.Range("D7").Value  = "6036f9b3b1842f8b659b18c7"
$ws.Range("E7").Value  = "Kellie"
$ws.Range("G7").Value  = 5.499920003737663
$ws.Range("H7").Value  = "White"

$ws.Range("G8").Value  = 5.496086788842061

$ws.Range("C9").Value  = 33
$ws.Range("D9").Value  = "60cb36ee9f58331a33cf5506"
$ws.Range("E9").Value  = "Shaniek"
$ws.Range("G9").Value  = 5.062422754775289
$ws.Range("H9").Value  = "Black or African American"

$ws.Range("C10").Value = 35
$ws.Range("D10").Value = "6077db0613ce87b4a62a78f9"
$ws.Range("E10").Value = "Lori"
$ws.Range("G10").Value = 4.244814854093466

$ws.Range("C11").Value = 34
$ws.Range("D11").Value = "5e96194b0a9fe909389e9f7b"
$ws.Range("E11").Value = "Tina"
$ws.Range("G11").Value = 4.079969157910064

$ws.Range("G12").Value = 2.218831050136576
$ws.Range("G13").Value = 1.453411958882284

# --- Male group (rows 14-25) ---
$ws.Range("G14").Value = 14.36820170540361
$ws.Range("G15").Value = 13.10445131131576
$ws.Range("G16").Value = 8.165085299637123
$ws.Range("G17").Value = 7.471687554102642
$ws.Range("G18").Value = 6.394477316427598
$ws.Range("G19").Value = 6.210208915507454
$ws.Range("G20").Value = 5.464785281559637

$ws.Range("C21").Value = 33
$ws.Range("D21").Value = "60b322994d0b901954690036"
$ws.Range("E21").Value = "Brennan"
$ws.Range("G21").Value = 5.311456126118004
$ws.Range("H21").Value = "White"

$ws.Range("C22").Value = 32
$ws.Range("D22").Value = "60bf9943e4e04642d4634ecc"
$ws.Range("E22").Value = "Jamarii"
$ws.Range("G22").Value = 5.011805588421218
$ws.Range("H22").Value = "Black or African American"

$ws.Range("G23").Value = 3.497412773125043
$ws.Range("G24").Value = 1.210064642988239
$ws.Range("G25").Value = 0.4333488792121737
